# Refresh cryptocurrency Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.001.85"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").Value = "2.630.37"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  -0.10%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "513.82"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "144.23"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("E7").Value = "  -0.52%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.570"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").Value = "2.656.59"
$ws.Range("E9").Value = "  +1.58%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.34"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.14%  "

$ws.Range("E11").Value = "  +1.62%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.338"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "3.087.94"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "58.954.38"
$ws.Range("E15").Value = "  -1.58%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.10"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "2.646.66"
$ws.Range("E18").Value = "  +0.78%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.55"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.24%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "344.51"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.45%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.38"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.11"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "61.11"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.421"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").Value = "2.739.97"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  -0.81%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.161"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").Value = "0.0₃0807"
$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("E30").Value = "  +2.57%  "

$ws.Range("E31").Value = "  -0.40%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.45"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +8.97%  "

$ws.Range("E33").Value = "  +0.06%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "18.89"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.54%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "150.11"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  +12.81%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.03"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.47%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.66%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.855"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.50"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("E42").Value = "  +0.15%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "281.34"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.40%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "

$ws.Range("E45").Value = "  -0.39%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0986"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "19.54"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.26%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0537"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "10.28"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  -0.56%  "

$ws.Range("D51").Value = "1.973.94"
$ws.Range("E51").Value = "  +1.27%  "
